# Auto-generated edit script: apply cell-value updates per the commit diff.
# Each sheet's numeric cells (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) are refreshed to new market-data snapshot values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 789.6
$ws.Range("I18").Value = 598.25
$ws.Range("J18").Value = 1555
$ws.Range("K18").Value = 598.25
$ws.Range("L18").Value = 1555
$ws.Range("M18").Value = -314.25
$ws.Range("N18").Value = -2123

$ws.Range("H74").Value = 2896.1482
$ws.Range("I74").Value = 2615.8333
$ws.Range("J74").Value = 3120.4
$ws.Range("K74").Value = 2615.8333
$ws.Range("L74").Value = 3120.4
$ws.Range("M74").Value = -1679.8333
$ws.Range("N74").Value = -4992.4

$ws.Range("H77").Value = 2896.1482
$ws.Range("I77").Value = 2615.8333
$ws.Range("J77").Value = 3120.4
$ws.Range("K77").Value = 13079.1665
$ws.Range("L77").Value = 15602
$ws.Range("M77").Value = -8399.166499999999
$ws.Range("N77").Value = -24962

$ws.Range("H88").Value = 21743.902
$ws.Range("I88").Value = 57655.7
$ws.Range("J88").Value = 4643.048
$ws.Range("K88").Value = 57655.7
$ws.Range("L88").Value = 4643.048
$ws.Range("M88").Value = -57249.7
$ws.Range("N88").Value = -5455.048

$ws.Range("H91").Value = 21743.902
$ws.Range("I91").Value = 57655.7
$ws.Range("J91").Value = 4643.048
$ws.Range("K91").Value = 57655.7
$ws.Range("L91").Value = 4643.048
$ws.Range("M91").Value = -56251.7
$ws.Range("N91").Value = -7451.048

$ws.Range("H132").Value = 4754.3335
$ws.Range("I132").Value = 1343.9286
$ws.Range("J132").Value = 52500
$ws.Range("K132").Value = 4031.7858
$ws.Range("L132").Value = 157500
$ws.Range("M132").Value = -1501.7858
$ws.Range("N132").Value = -162560


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11364675
$ws.Range("I45").Value = 12988057
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 12988057
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -12987680
$ws.Range("N45").Value = -1754

$ws.Range("H110").Value = 1345.1428
$ws.Range("I110").Value = 1500
$ws.Range("J110").Value = 1319.3334
$ws.Range("K110").Value = 1500
$ws.Range("L110").Value = 1319.3334
$ws.Range("M110").Value = 545
$ws.Range("N110").Value = -5409.3334


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 13000
$ws.Range("J50").Value = 13000
$ws.Range("L50").Value = 13000
$ws.Range("N50").Value = -14250

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("N51").Value = 0

$ws.Range("H59").Value = 16909.092
$ws.Range("J59").Value = 16909.092
$ws.Range("L59").Value = 16909.092
$ws.Range("N59").Value = -19199.092

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("N60").Value = 0

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("N61").Value = 0

$ws.Range("H68").Value = 581429.1
$ws.Range("J68").Value = 581429.1
$ws.Range("L68").Value = 581429.1
$ws.Range("N68").Value = -582927.1

$ws.Range("H71").Value = 581429.1
$ws.Range("J71").Value = 581429.1
$ws.Range("L71").Value = 1744287.3
$ws.Range("N71").Value = -1751775.3

$ws.Range("H74").Value = 23333.334
$ws.Range("J74").Value = 23333.334
$ws.Range("L74").Value = 23333.334
$ws.Range("N74").Value = -25081.334

$ws.Range("H77").Value = 23333.334
$ws.Range("J77").Value = 23333.334
$ws.Range("L77").Value = 70000.00199999999
$ws.Range("N77").Value = -78736.00199999999

$ws.Range("H99").Value = 3235.7144
$ws.Range("I99").Value = 3540
$ws.Range("J99").Value = 2475
$ws.Range("K99").Value = 3540
$ws.Range("L99").Value = 2475
$ws.Range("M99").Value = -2042
$ws.Range("N99").Value = -5471

$ws.Range("H126").Value = 3235.7144
$ws.Range("I126").Value = 3540
$ws.Range("J126").Value = 2475
$ws.Range("K126").Value = 10620
$ws.Range("L126").Value = 7425
$ws.Range("M126").Value = -8150
$ws.Range("N126").Value = -12365

$ws.Range("H141").Value = 47573.168
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 50079.816
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 50079.816
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -60439.816


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 998.6667
$ws.Range("I5").Value = 656.4
$ws.Range("K5").Value = 1969.2
$ws.Range("M5").Value = -1857.2

$ws.Range("H87").Value = 6967.5557
$ws.Range("I87").Value = 5618.8335
$ws.Range("J87").Value = 9665
$ws.Range("K87").Value = 16856.5005
$ws.Range("L87").Value = 28995
$ws.Range("M87").Value = -15608.5005
$ws.Range("N87").Value = -31491

$ws.Range("H90").Value = 6967.5557
$ws.Range("I90").Value = 5618.8335
$ws.Range("J90").Value = 9665
$ws.Range("K90").Value = 50569.5015
$ws.Range("L90").Value = 86985
$ws.Range("M90").Value = -44329.5015
$ws.Range("N90").Value = -99465

$ws.Range("H107").Value = 855158.1
$ws.Range("I107").Value = 467.83334
$ws.Range("J107").Value = 1010556.4
$ws.Range("K107").Value = 1403.50002
$ws.Range("L107").Value = 3031669.2
$ws.Range("M107").Value = 516.4999800000001
$ws.Range("N107").Value = -3035509.2

$ws.Range("H131").Value = 966.76
$ws.Range("I131").Value = 556.8333
$ws.Range("J131").Value = 992.92554
$ws.Range("K131").Value = 1670.4999
$ws.Range("L131").Value = 2978.77662
$ws.Range("M131").Value = 3369.5001
$ws.Range("N131").Value = -13058.77662

$ws.Range("H135").Value = 998.6667
$ws.Range("I135").Value = 656.4
$ws.Range("K135").Value = 5907.599999999999
$ws.Range("M135").Value = -3372.599999999999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2400
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -4396

$ws.Range("H83").Value = 2400
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -21984

$ws.Range("H102").Value = 4630594
$ws.Range("I102").Value = 5556337
$ws.Range("J102").Value = 1878.5
$ws.Range("K102").Value = 5556337
$ws.Range("L102").Value = 1878.5
$ws.Range("M102").Value = -5554715
$ws.Range("N102").Value = -5122.5

